$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append year data rows 199-210 (A = 197..208, B = tiny residual/normalized
# floats, "year data appended 1->12" / normalized ML X data).
# Scientific-notation literals aren't parsed by this shell, so values are
# written out in plain decimal form (identical underlying double value,
# just without the exponent shorthand).
$data = @(
    @(199, 197, 0.0000000000000001249000902703301),
    @(200, 198, -0.00000000000000004037174635000569),
    @(201, 199, -0.000000000000000005551115123125783),
    @(202, 200, -0.00000000000000002775557561562891),
    @(203, 201, -0.000000000000000135308431126191),
    @(204, 202, 0.0000000000000000515460690004537),
    @(205, 203, 0.00000000000000001156482317317871),
    @(206, 204, 0.0000000000000000791033905045424),
    @(207, 205, 0.00000000000000005551115123125783),
    @(208, 206, 0.00000000000000003700743415417188),
    @(209, 207, 0.00000000000000008326672684688674),
    @(210, 208, 0)
)

# Column A (rows 2-198) carries bold/centered/bordered formatting. Reuse it
# for the newly appended rows via copy/paste-format instead of rebuilding
# the style piecewise (which would mint extra, unused cellXf records).
$formatSource = $ws.Range("A198")

foreach ($entry in $data) {
    $row = $entry[0]
    $aVal = $entry[1]
    $bVal = $entry[2]

    $aCell = $ws.Range("A" + $row)
    $aCell.Value = $aVal
    $formatSource.Copy()
    $aCell.PasteSpecial(-4122)

    $ws.Range("B" + $row).Value = $bVal
}

$excel.CutCopyMode = $false
